# "adding points on map" - fixing mistyped city names so they match the
# points being geocoded/plotted on the map, and widening the City column
# so the corrected names are fully visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mistyped city names in column B (City)
$ws.Range("B74").Value = "Eastham"
$ws.Range("B76").Value = "Hingham"
$ws.Range("B93").Value = "Sudbury"
$ws.Range("B108").Value = "Billerica"
$ws.Range("B112").Value = "Medfield"

# Widen column B so the (now longer/corrected) city names display fully
$ws.Columns("B").ColumnWidth = 33.83

# Move the selection down to where the edits were made
$ws.Range("C113").Select()
